$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Rename "Population" -> "population"
$ws.Range("A3").Value = "population"

# Clear the whole Density column (C) - it's being removed as a column
# and replaced by a "density" row instead.
$ws.Range("C1:C3").Clear()

# Add new row 4: density label + value (moved from old C2/C3)
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 2517.176373984591
